$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Swap the two student records (row 2 <-> row 3) and update the
# "Reprobadas" (G) value for the student who now occupies row 3.

$ws.Range("A2").Value = 20330051920332
$ws.Range("B2").Value = "RODRIGUEZ"
$ws.Range("C2").Value = "GUERRA"
$ws.Range("D2").Value = "JAVIER ANTONIO"

$ws.Range("A3").Value = 20330051920323
$ws.Range("B3").Value = "GARCIA"
$ws.Range("C3").Value = "LEON"
$ws.Range("D3").Value = "JESUS SAMUEL"
$ws.Range("G3").Value = 1

$wb.Save()
